# "minor changes to sample data"
# The sample's anonymized patient codes (columns A "patient_id1" and B
# "patient_id2") are refreshed from the placeholder p1..p49 series to the
# real de-identified study codes (du0xx / iu-xx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$patientIds = @{
    2  = "du029"
    3  = "du032"
    4  = "du033"
    5  = "du036"
    6  = "du037"
    7  = "du038"
    8  = "du041"
    9  = "du042"
    10 = "du045"
    11 = "du046"
    12 = "du047"
    13 = "du048"
    14 = "du049"
    15 = "du051"
    16 = "du052"
    17 = "du054"
    18 = "du055"
    19 = "du057"
    20 = "iu-29"
    21 = "iu-31"
    22 = "iu-33"
    23 = "iu-34"
    24 = "iu-35"
    25 = "iu-38"
    26 = "iu-39"
    27 = "iu-40"
    28 = "iu-41"
    29 = "iu-42"
    30 = "iu-45"
    31 = "iu-47"
    32 = "iu-48"
    33 = "iu-49"
    34 = "iu-50"
    35 = "iu-55"
    36 = "iu-56"
    37 = "iu-57"
    38 = "iu-56"
    39 = "iu-57"
    40 = "iu-33"
    41 = "iu-34"
    42 = "iu-35"
    43 = "iu-38"
    44 = "iu-39"
    45 = "iu-40"
    46 = "iu-41"
    47 = "iu-42"
    48 = "iu-45"
    49 = "iu-47"
    50 = "iu-48"
}

foreach ($row in $patientIds.Keys) {
    $id = $patientIds[$row]
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $id
}

# Author's last selection before saving moved from L4 to C6.
$ws.Range("C6").Select() | Out-Null
